# Get the active workbook/worksheet (Backlog sheet is the active sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Id=4): task state changed from "Committed" to "On hold"
$ws.Range("B5").Value = "On hold"

# New row 19: a new backlog item
$ws.Range("A19").Value = 19
$ws.Range("B19").Value = "Approved"
$ws.Range("C19").Value = "Extract methods to helper class and unit test"

# Update the active selection to reflect where the user left off
$ws.Range("D27").Select() | Out-Null
